# Add 2 new scene quest rows (42010013 / 42010014) to Sheet1, right
# below the existing data, and grow the "表3" table to cover them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 22 - 丛林 (trees)
$ws.Cells.Item(22, 1).Value  = 42010013
$ws.Cells.Item(22, 2).Value  = "丛林"
$ws.Cells.Item(22, 3).Value  = 0
$ws.Cells.Item(22, 4).Value  = "trees"
$ws.Cells.Item(22, 5).Value  = "trees"
$ws.Cells.Item(22, 6).Value  = "trees"
$ws.Cells.Item(22, 8).Value  = 43000007
$ws.Cells.Item(22, 9).Value  = "mini"
$ws.Cells.Item(22, 10).Value = 50
$ws.Cells.Item(22, 11).Value = 150
$ws.Cells.Item(22, 19).Value = 100

# Row 23 - 墓地 (grave)
$ws.Cells.Item(23, 1).Value  = 42010014
$ws.Cells.Item(23, 2).Value  = "墓地"
$ws.Cells.Item(23, 3).Value  = 0
$ws.Cells.Item(23, 4).Value  = "grave"
$ws.Cells.Item(23, 5).Value  = "grave"
$ws.Cells.Item(23, 6).Value  = "grave"
$ws.Cells.Item(23, 8).Value  = 43000004
$ws.Cells.Item(23, 9).Value  = "oneline"
$ws.Cells.Item(23, 10).Value = 100
$ws.Cells.Item(23, 14).Value = 100
$ws.Cells.Item(23, 16).Value = 23000101
$ws.Cells.Item(23, 19).Value = 100
$ws.Cells.Item(23, 20).Value = 100

# Grow the worksheet table ("表3") so the two new rows join it.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A3:X23"))

# Move the active selection to match where the author ended up.
$ws.Range("S23").Select() | Out-Null
